$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "year"
$ws.Range("B1").Value = "objectName"
$ws.Range("C1").Value = "activity"
$ws.Range("D1").Value = "location"
$ws.Range("E1").Value = "no2"
$ws.Range("F1").Value = "so2"
$ws.Range("G1").Value = "co"
$ws.Range("H1").Value = "microparts"
$ws.Range("I1").Value = "summary"

$ws.Range("J1").Select()
